$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the description for the Throttling error code (row 36, A36=412034) to the
# more detailed wording (this is the actual content gap being filled in).
$ws.Range("B36").Value = "Throttling/Throttled - The connection is currently being throttled because the rate at which requests have been submitted for this action exceeds the limit for your account."

# Column B is sized to best-fit the longest description, so re-fit it now that the
# text in B36 is longer.
$ws.Columns.Item(2).AutoFit() | Out-Null

# Restore the view state (scroll position / active selection) left on the sheet
# after the edit.
$ws.Range("C28").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
